$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1069
    3  = 3
    5  = 3066
    7  = 2297
    11 = 1103
    15 = 474
    18 = 12
    19 = 9
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
